$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the C column dates to the new academic-year schedule ---
$ws.Range("C2").Value = 45891.427083333336
$ws.Range("C3").Value = 45894.46875
$ws.Range("C4").Value = 45898.427083333336
$ws.Range("C5").Value = 45901.46875
$ws.Range("C6").Value = 45905.427083333336
$ws.Range("C7").Value = 45908.46875
$ws.Range("C8").Value = 45912.427083333336
$ws.Range("C9").Value = 45915.46875
$ws.Range("C10").Value = 45919.427083333336
$ws.Range("C11").Value = 45922.46875
$ws.Range("C12").Value = 45926.427083333336
$ws.Range("C13").Value = 45929.46875
$ws.Range("C14").Value = 45933.427083333336
$ws.Range("C15").Value = 45936.46875
$ws.Range("C16").Value = 45940.427083333336
$ws.Range("C17").Value = 45950.46875
$ws.Range("C18").Value = 45954.427083333336
$ws.Range("C19").Value = 45957.46875
$ws.Range("C20").Value = 45961.427083333336
$ws.Range("C21").Value = 45968.427083333336
$ws.Range("C22").Value = 45971.46875
$ws.Range("C23").Value = 45975.427083333336
$ws.Range("C24").Value = 45978.46875
$ws.Range("C25").Value = 45982.427083333336
$ws.Range("C26").Value = 45985.46875
$ws.Range("C27").Value = 45989.427083333336
$ws.Range("C28").Value = 45992.46875
$ws.Range("C29").Value = 45996.427083333336
$ws.Range("C30").Value = 45999.46875
$ws.Range("C31").Value = 46003.427083333336

# --- Row height adjustments that followed the text re-wrap ---
$ws.Rows.Item(3).RowHeight = 109.2
$ws.Rows.Item(4).RowHeight = 78
$ws.Rows.Item(5).RowHeight = 46.8
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 62.4
$ws.Rows.Item(10).RowHeight = 78
$ws.Rows.Item(14).RowHeight = 31.2
$ws.Rows.Item(15).RowHeight = 62.4
$ws.Rows.Item(16).RowHeight = 109.2

# --- C5 loses its bottom/top border formatting (new un-bordered date style) ---
$ws.Range("C5").Borders.Item(8).LineStyle = -4142

# --- Append a new (still blank) row 32, carrying forward C31's date formatting ---
$ws.Range("C31").Copy()
$ws.Range("C32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the active selection to reflect the next empty row ---
$ws.Range("D33").Select() | Out-Null
